$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 157, pushing the existing rows 157-158 down to 158-159.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new weekly record.
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 45121
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = 100112031
$ws.Cells.Item(157, 7).Value = "Poroto verde"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 45
$ws.Cells.Item(157, 11).Value = 30000
$ws.Cells.Item(157, 12).Value = 30000
$ws.Cells.Item(157, 13).Value = 30000
$ws.Cells.Item(157, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(157, 15).Value = "Perú"
$ws.Cells.Item(157, 16).Value = 1200
$ws.Cells.Item(157, 17).Value = 25
$ws.Cells.Item(157, 18).Value = "Hortaliza"
